# Checkpoint before assistant change: Fix calculation of net profit by
# correctly subtracting expenses.
#
# This reshapes the header (and data) layout of the Products, Sales and
# Expenses sheets: an "id" column is introduced / relocated to the front,
# columns are reordered to match the application's canonical schema, a
# "receipt_number" column is added to Expenses, and the now-stale sample
# data row (row 2) is removed from each of these three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Products" (sheet1)
# ---------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.UsedRange.Clear()

$productHeaders = @(
    "id",
    "name",
    "description",
    "price",
    "cost_price",
    "category",
    "stock",
    "min_stock",
    "supplier",
    "sku",
    "created_date",
    "last_updated"
)
for ($i = 0; $i -lt $productHeaders.Length; $i++) {
    $wsProducts.Cells.Item(1, $i + 1).Value = $productHeaders[$i]
}

# ---------------------------------------------------------------------
# Sheet "Sales" (sheet2)
# ---------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item("Sales")
$wsSales.UsedRange.Clear()

$salesHeaders = @(
    "id",
    "product_id",
    "product_name",
    "quantity",
    "unit_price",
    "total_amount",
    "profit",
    "customer_name",
    "payment_method",
    "sale_date",
    "cashier",
    "notes"
)
for ($i = 0; $i -lt $salesHeaders.Length; $i++) {
    $wsSales.Cells.Item(1, $i + 1).Value = $salesHeaders[$i]
}

# ---------------------------------------------------------------------
# Sheet "Expenses" (sheet3)
# ---------------------------------------------------------------------
$wsExpenses = $wb.Worksheets.Item("Expenses")
$wsExpenses.UsedRange.Clear()

$expenseHeaders = @(
    "id",
    "category",
    "description",
    "amount",
    "payment_method",
    "vendor",
    "expense_date",
    "receipt_number",
    "notes"
)
for ($i = 0; $i -lt $expenseHeaders.Length; $i++) {
    $wsExpenses.Cells.Item(1, $i + 1).Value = $expenseHeaders[$i]
}

# "Goals" sheet (sheet4) is unchanged by this commit.
